$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update sheet (tab) name to reflect new "through" date
$ws.Name = "Through 2022-12-06"

# Update the December row label to reflect the new "through" date
$ws.Range("A13").Value = "December (through 12-06)"

# Update December figures (row 13) for years 2016-2022 (columns C-I)
$ws.Range("C13").Value = 19
$ws.Range("D13").Value = 23
$ws.Range("F13").Value = 6
$ws.Range("G13").Value = 31
$ws.Range("H13").Value = 49
$ws.Range("I13").Value = 24

# Update Total figures (row 14) for years 2016-2022 (columns C-I)
$ws.Range("C14").Value = 582
$ws.Range("D14").Value = 844
$ws.Range("F14").Value = 540
$ws.Range("G14").Value = 1295
$ws.Range("H14").Value = 1692
$ws.Range("I14").Value = 1539
